# Append the missing order row (row 5) to the inventory sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the new row as text first so that purely-numeric-looking
# values ("29", "234903243", "100") are stored as text, matching the
# style of the rest of the table, then restore the default "Normal"
# style so no extra formatting is left behind on the row.
$row = $ws.Range("A5:F5")
$row.NumberFormat = "@"

$ws.Range("A5").Value = "29"
$ws.Range("B5").Value = "234903243"
$ws.Range("C5").Value = "vintage lamp"
$ws.Range("D5").Value = "100"
$ws.Range("E5").Value = "bob"
$ws.Range("F5").Value = "dan"

$row.Style = "Normal"
